$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.186238
$ws.Range("H2").Value = 0.558714
$ws.Range("I2").Value = 0.05023668284714279
$ws.Range("J2").Value = 0.05023668284714279
$ws.Range("M2").Value = 4.093680666666667
$ws.Range("N2").Value = 12.281042
$ws.Range("O2").Value = 0.1610908176055751
$ws.Range("P2").Value = 0.161090817605575
$ws.Range("Q2").Value = 0.7623988999986667
$ws.Range("R2").Value = 6.861590099988
$ws.Range("S2").Value = 0.008092668313638201
$ws.Range("T2").Value = 0.0080926683136382

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.186238
$ws.Range("H3").Value = 0.558714
$ws.Range("I3").Value = 0.05023668284714279
$ws.Range("J3").Value = 0.05023668284714279
$ws.Range("O3").Value = 0.5606512265211691
$ws.Range("P3").Value = 0.5606512265211691
$ws.Range("Q3").Value = 2.653409329818
$ws.Range("R3").Value = 23.880683968362
$ws.Range("S3").Value = 0.02816525785460558
$ws.Range("T3").Value = 0.02816525785460558

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.186238
$ws.Range("H4").Value = 0.558714
$ws.Range("I4").Value = 0.05023668284714279
$ws.Range("J4").Value = 0.05023668284714279
$ws.Range("M4").Value = 7.071161666666666
$ws.Range("N4").Value = 21.213485
$ws.Range("O4").Value = 0.2782579558732559
$ws.Range("P4").Value = 0.2782579558732559
$ws.Range("Q4").Value = 1.316919006476667
$ws.Range("R4").Value = 11.85227105829
$ws.Range("S4").Value = 0.01397875667889901
$ws.Range("T4").Value = 0.01397875667889901

# Row 5
$ws.Range("I5").Value = 0.659992587420158
$ws.Range("J5").Value = 0.6599925874201579
$ws.Range("M5").Value = 4.093680666666667
$ws.Range("N5").Value = 12.281042
$ws.Range("O5").Value = 0.1610908176055751
$ws.Range("P5").Value = 0.161090817605575
$ws.Range("Q5").Value = 10.01613948491467
$ws.Range("R5").Value = 90.145255364232
$ws.Range("S5").Value = 0.1063187455211322
$ws.Range("T5").Value = 0.1063187455211322

# Row 6
$ws.Range("I6").Value = 0.659992587420158
$ws.Range("J6").Value = 0.6599925874201579
$ws.Range("O6").Value = 0.5606512265211691
$ws.Range("P6").Value = 0.5606512265211691
$ws.Range("S6").Value = 0.3700256536319915
$ws.Range("T6").Value = 0.3700256536319914

# Row 7
$ws.Range("I7").Value = 0.659992587420158
$ws.Range("J7").Value = 0.6599925874201579
$ws.Range("M7").Value = 7.071161666666666
$ws.Range("N7").Value = 21.213485
$ws.Range("O7").Value = 0.2782579558732559
$ws.Range("P7").Value = 0.2782579558732559
$ws.Range("Q7").Value = 17.30123752700667
$ws.Range("R7").Value = 155.71113774306
$ws.Range("S7").Value = 0.1836481882670343
$ws.Range("T7").Value = 0.1836481882670342

# Row 8
$ws.Range("G8").Value = 1.074241333333333
$ws.Range("H8").Value = 3.222724
$ws.Range("I8").Value = 0.2897707297326994
$ws.Range("J8").Value = 0.2897707297326994
$ws.Range("M8").Value = 4.093680666666667
$ws.Range("N8").Value = 12.281042
$ws.Range("O8").Value = 0.1610908176055751
$ws.Range("P8").Value = 0.161090817605575
$ws.Range("Q8").Value = 4.397600977600889
$ws.Range("R8").Value = 39.578408798408
$ws.Range("S8").Value = 0.04667940377080466
$ws.Range("T8").Value = 0.04667940377080465

# Row 9
$ws.Range("G9").Value = 1.074241333333333
$ws.Range("H9").Value = 3.222724
$ws.Range("I9").Value = 0.2897707297326994
$ws.Range("J9").Value = 0.2897707297326994
$ws.Range("O9").Value = 0.5606512265211691
$ws.Range("P9").Value = 0.5606512265211691
$ws.Range("Q9").Value = 15.305157789188
$ws.Range("R9").Value = 137.746420102692
$ws.Range("S9").Value = 0.1624603150345721
$ws.Range("T9").Value = 0.1624603150345721

# Row 10
$ws.Range("G10").Value = 1.074241333333333
$ws.Range("H10").Value = 3.222724
$ws.Range("I10").Value = 0.2897707297326994
$ws.Range("J10").Value = 0.2897707297326994
$ws.Range("M10").Value = 7.071161666666666
$ws.Range("N10").Value = 21.213485
$ws.Range("O10").Value = 0.2782579558732559
$ws.Range("P10").Value = 0.2782579558732559
$ws.Range("Q10").Value = 7.596134137015556
$ws.Range("R10").Value = 68.36520723314
$ws.Range("S10").Value = 0.08063101092732261
$ws.Range("T10").Value = 0.08063101092732261

